$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so exact formatting (trailing zeros, etc.) is preserved
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply updated values
# Row 2
$ws.Range("D2").Value = "29.923.25"
# Row 3
$ws.Range("D3").Value = "1.925.21"
$ws.Range("E3").Value = "  +1.79%  "
# Row 4
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.14%  "
# Row 5
$ws.Range("D5").Value = "320.60"
$ws.Range("E5").Value = "  -0.70%  "
# Row 6
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.11%  "
# Row 7
$ws.Range("D7").Value = "0.5059"
$ws.Range("E7").Value = "  -2.19%  "
# Row 8
$ws.Range("D8").Value = "0.4056"
$ws.Range("E8").Value = "  +1.12%  "
# Row 9
$ws.Range("D9").Value = "0.08349"
$ws.Range("E9").Value = "  -0.54%  "
# Row 10
$ws.Range("D10").Value = "42.40"
$ws.Range("E10").Value = "  -0.68%  "
# Row 11
$ws.Range("E11").Value = "  -0.70%  "
# Row 12
$ws.Range("D12").Value = "23.89"
$ws.Range("E12").Value = "  +3.80%  "
# Row 13
$ws.Range("D13").Value = "1.938.08"
$ws.Range("E13").Value = "  +2.34%  "
# Row 15
$ws.Range("D15").Value = "7.255"
$ws.Range("E15").Value = "  -0.69%  "
# Row 16
$ws.Range("D16").Value = "0.9986"
$ws.Range("E16").Value = "  -0.38%  "
# Row 17
$ws.Range("D17").Value = "92.29"
$ws.Range("E17").Value = "  -1.95%  "
# Row 18
$ws.Range("E18").Value = "  -0.78%  "
# Row 19
$ws.Range("D19").Value = "0.06518"
$ws.Range("E19").Value = "  -1.81%  "
# Row 20
$ws.Range("D20").Value = "18.28"
$ws.Range("E20").Value = "  +0.50%  "
# Row 21
$ws.Range("D21").Value = "0.9994"
$ws.Range("E21").Value = "  -0.06%  "
# Row 22
$ws.Range("D22").Value = "5.959"
$ws.Range("E22").Value = "  +0.33%  "
# Row 23
$ws.Range("D23").Value = "29.990.13"
$ws.Range("E23").Value = "  -0.82%  "
# Row 24
$ws.Range("D24").Value = "11.34"
$ws.Range("E24").Value = "  +0.67%  "
# Row 25
$ws.Range("D25").Value = "2.191"
$ws.Range("E25").Value = "  -1.73%  "
# Row 26
$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "2.151.01"
$ws.Range("E26").Value = "  +1.81%  "
# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "22.15"
$ws.Range("E27").Value = "  +2.79%  "
# Row 28
$ws.Range("D28").Value = "162.16"
# Row 29
$ws.Range("D29").Value = "2.338"
$ws.Range("E29").Value = "  +0.40%  "
# Row 30
$ws.Range("D30").Value = "129.04"
$ws.Range("E30").Value = "  +0.05%  "
# Row 31
$ws.Range("D31").Value = "1.132"
$ws.Range("E31").Value = "  +4.32%  "
# Row 32
$ws.Range("E32").Value = "  -1.34%  "
# Row 33
$ws.Range("D33").Value = "5.964"
$ws.Range("E33").Value = "  -1.97%  "
# Row 34
$ws.Range("D34").Value = "3.786"
$ws.Range("E34").Value = "  +1.00%  "
# Row 35
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "5.415"
$ws.Range("E35").Value = "  +1.72%  "
# Row 36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "0.02451"
$ws.Range("E36").Value = "  -1.49%  "
# Row 37
$ws.Range("D37").Value = "0.06431"
$ws.Range("E37").Value = "  -1.35%  "
# Row 38
$ws.Range("D38").Value = "0.2159"
$ws.Range("E38").Value = "  -1.48%  "
# Row 39
$ws.Range("D39").Value = "0.6572"
$ws.Range("E39").Value = "  +1.35%  "
# Row 40
$ws.Range("D40").Value = "8.776"
$ws.Range("E40").Value = "  -0.22%  "
# Row 41
$ws.Range("D41").Value = "1.198"
$ws.Range("E41").Value = "  -1.57%  "
# Row 43
$ws.Range("D43").Value = "1.214"
$ws.Range("E43").Value = "  -1.16%  "
# Row 44
$ws.Range("D44").Value = "2.243"
$ws.Range("E44").Value = "  +9.59%  "
# Row 45
$ws.Range("D45").Value = "13.40"
$ws.Range("E45").Value = "  +1.53%  "
# Row 46
$ws.Range("D46").Value = "0.6101"
$ws.Range("E46").Value = "  +0.44%  "
# Row 47
$ws.Range("D47").Value = "3.596"
$ws.Range("E47").Value = "  -2.39%  "
# Row 48
$ws.Range("D48").Value = "1.211"
$ws.Range("E48").Value = "  -1.88%  "
# Row 49
$ws.Range("D49").Value = "122.14"
$ws.Range("E49").Value = "  -1.75%  "
# Row 50
$ws.Range("D50").Value = "79.24"
$ws.Range("E50").Value = "  +0.43%  "
# Row 51
$ws.Range("E51").Value = "  -2.81%  "
